# Split the old "Terms Typically Offered" column (D) into separate
# requirement columns by inserting three new columns before it, then
# populate the new columns (D:F) with the appropriate header/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank columns at D:F, pushing the existing D column
# (Terms Typically Offered) to G.
$ws.Columns("D:F").Insert()

# New header row.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# New data rows - all "NA" for the sample data.
$ws.Range("D2:F5").Value = "NA"
